$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "ShowName" column in U (column 21), right after the existing
# "Icon" column in T (column 20).
$ws.Cells.Item(1, 21).Value = "ShowName"

# Column T currently stores icon paths with a ".png" suffix (e.g.
# "UI/SteampunkUI/resource/icons/img_equip.png"); strip the extension, and
# populate the new column U with "ShowName_<n>" for each data row.
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $iconCell = $ws.Cells.Item($r, 20)
    $iconVal = $iconCell.Value2
    if ($iconVal -ne $null -and $iconVal -like "*.png") {
        $iconCell.Value = $iconVal.Substring(0, $iconVal.Length - 4)
    }

    $ws.Cells.Item($r, 21).Value = "ShowName_" + ($r - 1)
}

# Match the author's final selection/scroll position after the edit.
$ws.Range("W43").Select() | Out-Null

